# Insert a new data row at row 73 (shifts existing rows 73..112 down to 74..113)
# and populate it with the new weekly observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(73).Insert()

$ws.Cells.Item(73, 1).Value = 3
$ws.Cells.Item(73, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(73, 3).Value = "Coquimbo"
$ws.Cells.Item(73, 4).Value = 44466
$ws.Cells.Item(73, 5).Value = 5
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100101
$ws.Cells.Item(73, 8).Value = "Berries"
$ws.Cells.Item(73, 9).Value = 100112025
$ws.Cells.Item(73, 10).Value = "Frutilla"
$ws.Cells.Item(73, 11).Value = "Sin especificar"
$ws.Cells.Item(73, 12).Value = "Especial"
$ws.Cells.Item(73, 13).Value = 70
$ws.Cells.Item(73, 14).Value = 15000
$ws.Cells.Item(73, 15).Value = 15000
$ws.Cells.Item(73, 16).Value = 15000
$ws.Cells.Item(73, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(73, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(73, 19).Value = 2143
$ws.Cells.Item(73, 20).Value = 7
